# Weekly update: insert 4 new price records at the top of the data block
# (rows 1140-1143), pushing the existing rows 1140-1227 down to 1144-1231.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at row 1140 (entire rows shift down, same as
# right-click "Insert" on the row headers in Excel).
$ws.Rows("1140:1143").Insert()

# Common/static values shared by every record in this data set.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112021
$categoria = "Ají"
$clasificacion = "Hortaliza"

# New row 1140
$ws.Cells.Item(1140, 1).Value = $mercadoId
$ws.Cells.Item(1140, 2).Value = $mercado
$ws.Cells.Item(1140, 3).Value = $region
$ws.Cells.Item(1140, 4).Value = 45013
$ws.Cells.Item(1140, 5).Value = $codreg
$ws.Cells.Item(1140, 6).Value = $categoriaId
$ws.Cells.Item(1140, 7).Value = $categoria
$ws.Cells.Item(1140, 8).Value = "Americana (o)"
$ws.Cells.Item(1140, 9).Value = "Primera"
$ws.Cells.Item(1140, 10).Value = 80
$ws.Cells.Item(1140, 11).Value = 15000
$ws.Cells.Item(1140, 12).Value = 17000
$ws.Cells.Item(1140, 13).Value = 16125
$ws.Cells.Item(1140, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(1140, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(1140, 16).Value = 645
$ws.Cells.Item(1140, 17).Value = 25
$ws.Cells.Item(1140, 18).Value = $clasificacion

# New row 1141
$ws.Cells.Item(1141, 1).Value = $mercadoId
$ws.Cells.Item(1141, 2).Value = $mercado
$ws.Cells.Item(1141, 3).Value = $region
$ws.Cells.Item(1141, 4).Value = 45013
$ws.Cells.Item(1141, 5).Value = $codreg
$ws.Cells.Item(1141, 6).Value = $categoriaId
$ws.Cells.Item(1141, 7).Value = $categoria
$ws.Cells.Item(1141, 8).Value = "Americana (o)"
$ws.Cells.Item(1141, 9).Value = "Primera"
$ws.Cells.Item(1141, 10).Value = 400
$ws.Cells.Item(1141, 11).Value = 17000
$ws.Cells.Item(1141, 12).Value = 18000
$ws.Cells.Item(1141, 13).Value = 17425
$ws.Cells.Item(1141, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(1141, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1141, 16).Value = 697
$ws.Cells.Item(1141, 17).Value = 25
$ws.Cells.Item(1141, 18).Value = $clasificacion

# New row 1142
$ws.Cells.Item(1142, 1).Value = $mercadoId
$ws.Cells.Item(1142, 2).Value = $mercado
$ws.Cells.Item(1142, 3).Value = $region
$ws.Cells.Item(1142, 4).Value = 45013
$ws.Cells.Item(1142, 5).Value = $codreg
$ws.Cells.Item(1142, 6).Value = $categoriaId
$ws.Cells.Item(1142, 7).Value = $categoria
$ws.Cells.Item(1142, 8).Value = "Americana (o)"
$ws.Cells.Item(1142, 9).Value = "Segunda"
$ws.Cells.Item(1142, 10).Value = 120
$ws.Cells.Item(1142, 11).Value = 15000
$ws.Cells.Item(1142, 12).Value = 15000
$ws.Cells.Item(1142, 13).Value = 15000
$ws.Cells.Item(1142, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(1142, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1142, 16).Value = 600
$ws.Cells.Item(1142, 17).Value = 25
$ws.Cells.Item(1142, 18).Value = $clasificacion

# New row 1143
$ws.Cells.Item(1143, 1).Value = $mercadoId
$ws.Cells.Item(1143, 2).Value = $mercado
$ws.Cells.Item(1143, 3).Value = $region
$ws.Cells.Item(1143, 4).Value = 45013
$ws.Cells.Item(1143, 5).Value = $codreg
$ws.Cells.Item(1143, 6).Value = $categoriaId
$ws.Cells.Item(1143, 7).Value = $categoria
$ws.Cells.Item(1143, 8).Value = "Cacho cabra rojo"
$ws.Cells.Item(1143, 9).Value = "Primera"
$ws.Cells.Item(1143, 10).Value = 400
$ws.Cells.Item(1143, 11).Value = 18000
$ws.Cells.Item(1143, 12).Value = 19000
$ws.Cells.Item(1143, 13).Value = 18425
$ws.Cells.Item(1143, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(1143, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1143, 16).Value = 737
$ws.Cells.Item(1143, 17).Value = 25
$ws.Cells.Item(1143, 18).Value = $clasificacion
